$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 7051.8237
$ws.Range("I43").Value = 13459.875
$ws.Range("J43").Value = 1355.7778
$ws.Range("K43").Value = 13459.875
$ws.Range("L43").Value = 1355.7778
$ws.Range("M43").Value = -13390.875
$ws.Range("N43").Value = -1493.7778
$ws.Range("H51").Value = 1675.2307
$ws.Range("I51").Value = 1293.3334
$ws.Range("K51").Value = 1293.3334
$ws.Range("M51").Value = -809.3334
$ws.Range("H64").Value = 3793.75
$ws.Range("I64").Value = 3750
$ws.Range("J64").Value = 3925
$ws.Range("K64").Value = 3750
$ws.Range("L64").Value = 3925
$ws.Range("M64").Value = -3502
$ws.Range("N64").Value = -4421
$ws.Range("H67").Value = 3793.75
$ws.Range("I67").Value = 3750
$ws.Range("J67").Value = 3925
$ws.Range("K67").Value = 3750
$ws.Range("L67").Value = 3925
$ws.Range("M67").Value = -2892
$ws.Range("N67").Value = -5641
$ws.Range("H113").Value = 11306.25
$ws.Range("I113").Value = 21430.4
$ws.Range("J113").Value = 4074.7144
$ws.Range("K113").Value = 21430.4
$ws.Range("L113").Value = 4074.7144
$ws.Range("M113").Value = -18176.4
$ws.Range("N113").Value = -10582.7144
$ws.Range("H116").Value = 290677.16
$ws.Range("I116").Value = 3307.875
$ws.Range("J116").Value = 532672.3
$ws.Range("K116").Value = 3307.875
$ws.Range("L116").Value = 532672.3
$ws.Range("M116").Value = 134.125
$ws.Range("N116").Value = -539556.3
$ws.Range("H129").Value = 844.0345
$ws.Range("I129").Value = 421.30768
$ws.Range("J129").Value = 1187.5
$ws.Range("K129").Value = 1263.92304
$ws.Range("L129").Value = 3562.5
$ws.Range("M129").Value = 3736.07696
$ws.Range("N129").Value = -13562.5
$ws.Range("H132").Value = 43435736
$ws.Range("I132").Value = 53087396
$ws.Range("J132").Value = 3251
$ws.Range("K132").Value = 159262188
$ws.Range("L132").Value = 9753
$ws.Range("M132").Value = -159259658
$ws.Range("N132").Value = -14813

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2557.2375
$ws.Range("I32").Value = 2276.6082
$ws.Range("J32").Value = 6018.3335
$ws.Range("K32").Value = 2276.6082
$ws.Range("L32").Value = 6018.3335
$ws.Range("M32").Value = -1989.6082
$ws.Range("N32").Value = -6592.3335
$ws.Range("H61").Value = 5049.222
$ws.Range("I61").Value = 6129.9473
$ws.Range("J61").Value = 2482.5
$ws.Range("K61").Value = 6129.9473
$ws.Range("L61").Value = 2482.5
$ws.Range("M61").Value = -5917.9473
$ws.Range("N61").Value = -2906.5
$ws.Range("H74").Value = 5311.24
$ws.Range("I74").Value = 1140.0769
$ws.Range("J74").Value = 9830
$ws.Range("K74").Value = 1140.0769
$ws.Range("L74").Value = 9830
$ws.Range("M74").Value = -266.0769
$ws.Range("N74").Value = -11578
$ws.Range("H77").Value = 5311.24
$ws.Range("I77").Value = 1140.0769
$ws.Range("J77").Value = 9830
$ws.Range("K77").Value = 5700.3845
$ws.Range("L77").Value = 49150
$ws.Range("M77").Value = -1332.3845
$ws.Range("N77").Value = -57886
$ws.Range("H132").Value = 2779113.2
$ws.Range("I132").Value = 3572268.8
$ws.Range("J132").Value = 3068.8
$ws.Range("K132").Value = 10716806.4
$ws.Range("L132").Value = 9206.400000000001
$ws.Range("M132").Value = -10714276.4
$ws.Range("N132").Value = -14266.4
$ws.Range("H134").Value = 17300
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 17300
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 17300
$ws.Range("M134").ClearContents()  # was -5320, becomes blank (profit positive, omitted)
$ws.Range("N134").Value = -27440
$ws.Range("H136").Value = 5049.222
$ws.Range("I136").Value = 6129.9473
$ws.Range("J136").Value = 2482.5
$ws.Range("K136").Value = 18389.8419
$ws.Range("L136").Value = 7447.5
$ws.Range("M136").Value = -15839.8419
$ws.Range("N136").Value = -12547.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 18563340
$ws.Range("I80").Value = 183751.67
$ws.Range("J80").Value = 23814652
$ws.Range("K80").Value = 183751.67
$ws.Range("L80").Value = 23814652
$ws.Range("M80").Value = -182753.67
$ws.Range("N80").Value = -23816648
$ws.Range("H83").Value = 18563340
$ws.Range("I83").Value = 183751.67
$ws.Range("J83").Value = 23814652
$ws.Range("K83").Value = 918758.3500000001
$ws.Range("L83").Value = 119073260
$ws.Range("M83").Value = -913766.3500000001
$ws.Range("N83").Value = -119083244
$ws.Range("H134").Value = 11129259
$ws.Range("I134").Value = 12841091
$ws.Range("K134").Value = 38523273
$ws.Range("M134").Value = -38520738

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2473.2727
$ws.Range("I62").Value = 2171.4285
$ws.Range("J62").Value = 3001.5
$ws.Range("K62").Value = 2171.4285
$ws.Range("L62").Value = 3001.5
$ws.Range("M62").Value = -1547.4285
$ws.Range("N62").Value = -4249.5
$ws.Range("H65").Value = 2473.2727
$ws.Range("I65").Value = 2171.4285
$ws.Range("J65").Value = 3001.5
$ws.Range("K65").Value = 10857.1425
$ws.Range("L65").Value = 15007.5
$ws.Range("M65").Value = -7737.1425
$ws.Range("N65").Value = -21247.5
$ws.Range("H122").Value = 5330156
$ws.Range("I122").Value = 11906271
$ws.Range("J122").Value = 69264.266
$ws.Range("K122").Value = 35718813
$ws.Range("L122").Value = 207792.798
$ws.Range("M122").Value = -35716363
$ws.Range("N122").Value = -212692.798
$ws.Range("H132").Value = 9263356
$ws.Range("I132").Value = 12821326
$ws.Range("J132").Value = 12634.2
$ws.Range("K132").Value = 38463978
$ws.Range("L132").Value = 37902.60000000001
$ws.Range("M132").Value = -38461448
$ws.Range("N132").Value = -42962.60000000001
$ws.Range("H134").Value = 12020720
$ws.Range("I134").Value = 16668098
$ws.Range("J134").Value = 5683386.5
$ws.Range("K134").Value = 50004294
$ws.Range("L134").Value = 17050159.5
$ws.Range("M134").Value = -50001759
$ws.Range("N134").Value = -17055229.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1566.9756
$ws.Range("I126").Value = 1173.12
$ws.Range("J126").Value = 2182.375
$ws.Range("K126").Value = 3519.36
$ws.Range("L126").Value = 6547.125
$ws.Range("M126").Value = -1049.36
$ws.Range("N126").Value = -11487.125
$ws.Range("H132").Value = 55561148
$ws.Range("I132").Value = 90910430
$ws.Range("J132").Value = 12281.714
$ws.Range("K132").Value = 272731290
$ws.Range("L132").Value = 36845.142
$ws.Range("M132").Value = -272728760
$ws.Range("N132").Value = -41905.142

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2002.2222
$ws.Range("I40").Value = 2074.2856
$ws.Range("J40").Value = 1750
$ws.Range("K40").Value = 2074.2856
$ws.Range("L40").Value = 1750
$ws.Range("M40").Value = -1938.2856
$ws.Range("N40").Value = -2022
$ws.Range("H55").Value = 10327.6
$ws.Range("I55").Value = 567.3333
$ws.Range("J55").Value = 14510.571
$ws.Range("K55").Value = 567.3333
$ws.Range("L55").Value = 14510.571
$ws.Range("M55").Value = -394.3333
$ws.Range("N55").Value = -14856.571

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 19933392
$ws.Range("I136").Value = 14879276
$ws.Range("J136").Value = 26317538
$ws.Range("K136").Value = 44637828
$ws.Range("L136").Value = 78952614
$ws.Range("M136").Value = -44635278
$ws.Range("N136").Value = -78957714
